# Updated cryptos list on Thu Oct 10 13:30:22 UTC 2024 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCellValue($range, [string]$value) {
    # Force text storage so numeric-looking strings (e.g. with
    # trailing zeros or many leading zeros) are not reinterpreted
    # as numbers, then restore the default "Normal" style so no
    # stray number-format style is left behind on the cell.
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

# Row 2
$ws.Range('D2').Value = '61.184.63'
$ws.Range('E2').Value = '  -1.19%  '

# Row 3
$ws.Range('D3').Value = '2.417.82'
$ws.Range('E3').Value = '  -1.04%  '

# Row 4
$ws.Range('D4').Value = '0.996'
$ws.Range('E4').Value = '  -0.21%  '

# Row 5
$ws.Range('D5').Value = '568.58'
$ws.Range('E5').Value = '  -1.79%  '

# Row 6
$ws.Range('D6').Value = '139.98'
$ws.Range('E6').Value = '  -0.81%  '

# Row 7
$ws.Range('E7').Value = '  +0.01%  '

# Row 8
$ws.Range('D8').Value = '0.537'
$ws.Range('E8').Value = '  +1.36%  '

# Row 9
$ws.Range('D9').Value = '2.403.41'
$ws.Range('E9').Value = '  -1.45%  '

# Row 10
$ws.Range('D10').Value = '0.107'
$ws.Range('E10').Value = '  -1.95%  '

# Row 11
$ws.Range('E11').Value = '  -0.41%  '

# Row 12
$ws.Range('D12').Value = '5.07'
$ws.Range('E12').Value = '  -1.81%  '

# Row 13
$ws.Range('D13').Value = '0.337'
$ws.Range('E13').Value = '  -0.66%  '

# Row 14
$ws.Range('D14').Value = '26.11'
$ws.Range('E14').Value = '  +0.79%  '

# Row 15
$ws.Range('B15').Value = 'ShibaInu'
$ws.Range('C15').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
Set-TextCellValue $ws.Range('D15') '0.0000170'
$ws.Range('E15').Value = '  -1.34%  '

# Row 16
$ws.Range('B16').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C16').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D16').Value = '2.854.84'
$ws.Range('E16').Value = '  -1.33%  '

# Row 17
$ws.Range('D17').Value = '61.098.38'
$ws.Range('E17').Value = '  -1.27%  '

# Row 18
$ws.Range('D18').Value = '2.401.66'
$ws.Range('E18').Value = '  -1.62%  '

# Row 19
$ws.Range('D19').Value = '8.03'
$ws.Range('E19').Value = '  +11.51%  '

# Row 20
$ws.Range('D20').Value = '10.55'
$ws.Range('E20').Value = '  -0.35%  '

# Row 21
$ws.Range('D21').Value = '323.31'
$ws.Range('E21').Value = '  -0.51%  '

# Row 22
$ws.Range('D22').Value = '4.05'
$ws.Range('E22').Value = '  +0.08%  '

# Row 23
$ws.Range('E23').Value = '  +3.85%  '

# Row 24
$ws.Range('E24').Value = '  +0.09%  '

# Row 25
$ws.Range('D25').Value = '1.85'
$ws.Range('E25').Value = '  -2.47%  '

# Row 26
$ws.Range('D26').Value = '64.43'
$ws.Range('E26').Value = '  -0.94%  '

# Row 27
$ws.Range('D27').Value = '584.92'
$ws.Range('E27').Value = '  -0.36%  '

# Row 28
Set-TextCellValue $ws.Range('D28') '8.30'
$ws.Range('E28').Value = '  -9.27%  '

# Row 29
$ws.Range('D29').Value = '2.542.84'
$ws.Range('E29').Value = '  +0.35%  '

# Row 30
$ws.Range('D30').Value = '0.0₃0932'
$ws.Range('E30').Value = '  -0.82%  '

# Row 31
$ws.Range('D31').Value = '7.98'
$ws.Range('E31').Value = '  +1.35%  '

# Row 32
$ws.Range('D32').Value = '1.34'
$ws.Range('E32').Value = '  -3.10%  '

# Row 33
$ws.Range('D33').Value = '1.82'
$ws.Range('E33').Value = '  -2.94%  '

# Row 34
$ws.Range('E34').Value = '  -0.79%  '

# Row 35
$ws.Range('E35').Value = '  +0.02%  '

# Row 36
$ws.Range('B36').Value = 'ImmutableX'
$ws.Range('C36').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D36').Value = '1.41'
$ws.Range('E36').Value = '  +1.01%  '

# Row 37
$ws.Range('B37').Value = 'Monero'
$ws.Range('C37').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D37').Value = '151.55'
$ws.Range('E37').Value = '  -0.86%  '

# Row 38
$ws.Range('B38').Value = 'NEARProtocol'
$ws.Range('C38').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D38').Value = '4.59'
$ws.Range('E38').Value = '  -3.69%  '

# Row 39
$ws.Range('B39').Value = 'PolygonEcosystemToken'
$ws.Range('C39').Value = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
$ws.Range('D39').Value = '0.369'
$ws.Range('E39').Value = '  -1.20%  '

# Row 40
$ws.Range('D40').Value = '18.23'
$ws.Range('E40').Value = '  -0.57%  '

# Row 41
$ws.Range('D41').Value = '5.14'
$ws.Range('E41').Value = '  -0.74%  '

# Row 42
$ws.Range('D42').Value = '0.999'
$ws.Range('E42').Value = '  -0.03%  '

# Row 43
$ws.Range('B43').Value = 'Stacks'
$ws.Range('C43').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D43').Value = '1.67'
$ws.Range('E43').Value = '  -0.65%  '

# Row 44
$ws.Range('B44').Value = 'OKB'
$ws.Range('C44').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D44').Value = '41.19'
$ws.Range('E44').Value = '  -3.04%  '

# Row 45
$ws.Range('D45').Value = '2.41'
$ws.Range('E45').Value = '  +1.72%  '

# Row 46
$ws.Range('D46').Value = '0.0₆0295'
$ws.Range('E46').Value = '  +9.31%  '

# Row 47
Set-TextCellValue $ws.Range('D47') '143.10'
$ws.Range('E47').Value = '  +1.74%  '

# Row 48
$ws.Range('D48').Value = '3.53'
$ws.Range('E48').Value = '  -1.53%  '

# Row 49
$ws.Range('D49').Value = '0.587'
$ws.Range('E49').Value = '  -1.72%  '

# Row 50
$ws.Range('B50').Value = 'InjectiveProtocol'
$ws.Range('C50').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D50').Value = '19.36'
$ws.Range('E50').Value = '  -1.04%  '

# Row 51
$ws.Range('B51').Value = 'Hedera'
$ws.Range('C51').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D51').Value = '0.0501'
$ws.Range('E51').Value = '  -2.05%  '
